# Exercise 5.1 rewrite: array de-dup -> count-of-identical-run-lengths.
$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# 1) Problem statement paragraph (originally mentions "Integers").
$idx1 = Find-ParagraphIndex $d "Integers"
if ($idx1 -eq -1) { $idx1 = Find-ParagraphIndex $d "כתבו פעולה" }
$xmlPara1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve">כתבו פעולה שמקבלת כקלט </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>מערך</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve"> של </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>מספרים שלמים</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve">, ומחזירה </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve">מספר שלם אשר שווה למספר הרצפים של מספרים זהים במערך.  </w:t></w:r></w:p>'
$d.Paragraphs($idx1).Range.InsertXML($xmlPara1)

# 2) Example paragraph (originally mentions "לדוגמא").
$idx2 = Find-ParagraphIndex $d "לדוגמא"
$xmlPara2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="cs"/><w:b/><w:bCs/><w:rtl/></w:rPr><w:t xml:space="preserve">לדוגמא: </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve">על קלט </w:t></w:r><w:r><w:rPr><w:rtl/></w:rPr><w:t>[5, 5, 1, 3, 5, 3, 3, 3]</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve"> יוחזר </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>המספר</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>5</w:t></w:r></w:p>'
$d.Paragraphs($idx2).Range.InsertXML($xmlPara2)

# 3) Code signature paragraph (originally declares "countChanges").
$idx3 = Find-ParagraphIndex $d "countChanges"
if ($idx3 -eq -1) { $idx3 = Find-ParagraphIndex $d "public static int" }
$xmlPara3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="2B2B2B"/><w:bidi w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="A9B7C6"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-IL" w:eastAsia="en-IL"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="CC7832"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-IL" w:eastAsia="en-IL"/></w:rPr><w:t xml:space="preserve">public static int </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="FFC66D"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-IL" w:eastAsia="en-IL"/></w:rPr><w:t>countIdenticalSequences</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="A9B7C6"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-IL" w:eastAsia="en-IL"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="CC7832"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-IL" w:eastAsia="en-IL"/></w:rPr><w:t xml:space="preserve">int </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="A9B7C6"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-IL" w:eastAsia="en-IL"/></w:rPr><w:t>arr</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="A9B7C6"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-IL" w:eastAsia="en-IL"/></w:rPr><w:t>[])</w:t></w:r></w:p>'
$d.Paragraphs($idx3).Range.InsertXML($xmlPara3)

# 4) Blank paragraph right after the code block.
$idx4 = $idx3 + 1
$xmlPara4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rtl/><w:lang w:val="en-IL"/></w:rPr></w:pPr></w:p>'
$d.Paragraphs($idx4).Range.InsertXML($xmlPara4)

Write-Host "Updated paragraphs:" $idx1 $idx2 $idx3 $idx4
